$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.678.03"
$ws.Range("E2").Value = "  -0.41%  "

$ws.Range("D3").Value = "2.396.57"
$ws.Range("E3").Value = "  -0.87%  "

$ws.Range("E4").Value = "  +0.79%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "

$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("E8").Value = "  +2.17%  "

$ws.Range("D9").Value = "2.402.56"
$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("E10").Value = "  +0.79%  "

$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.97%  "

$ws.Range("E13").Value = "  +2.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("D17").Value = "60.287.32"
$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").Value = "2.400.25"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.67"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("E22").Value = "  +1.29%  "

$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "579.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.64%  "

$ws.Range("D30").Value = "0.0₃0935"
$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.11%  "

$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("E33").Value = "  -1.51%  "

$ws.Range("E34").Value = "  +0.29%  "

$ws.Range("E35").Value = "  -0.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.41%  "

$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("E43").Value = "  +6.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.18%  "

$ws.Range("D46").Value = "0.0₆0279"
$ws.Range("E46").Value = "  +4.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("E49").Value = "  +0.58%  "

$ws.Range("E50").Value = "  +1.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "

